# Add the "lostsoul" entity to the collision table: new row 10 and new column J,
# mirroring the existing pattern (collides with everything except "player").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "lostsoul" goes in header cell A10 and J1.
$ws.Range("A10").Value = "lostsoul"
$ws.Range("J1").Value = "lostsoul"

# New column J (rows 2-10): collision values, mirroring row 10.
$ws.Range("J2").Value = "▬"
$ws.Range("J3").Value = "x"
$ws.Range("J4").Value = "x"
$ws.Range("J5").Value = "x"
$ws.Range("J6").Value = "x"
$ws.Range("J7").Value = "x"
$ws.Range("J8").Value = "x"
$ws.Range("J9").Value = "x"
$ws.Range("J10").Value = "x"

# New row 10 (columns B-J): collision values.
$ws.Range("B10").Value = "▬"
$ws.Range("C10").Value = "x"
$ws.Range("D10").Value = "x"
$ws.Range("E10").Value = "x"
$ws.Range("F10").Value = "x"
$ws.Range("G10").Value = "x"
$ws.Range("H10").Value = "x"
$ws.Range("I10").Value = "x"

# Column J was narrower than the default (Excel auto-shrank it because
# "lostsoul" is shorter than "firefoebullet"/"waterfoebullet" etc.).
$ws.Columns.Item(10).ColumnWidth = 13

# Update selection to match the post-edit state.
$ws.Range("F10").Select()
